$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells we touch keep a Text format so values like "30.412.20",
# "0.000007621" or "0.1210" are not reinterpreted as numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.412.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.03"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9992"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4687"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.49%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06814"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "111.13"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +14.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.39"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07724"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.889.95"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.258"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6585"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.79%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.422.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007621"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.135.99"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9982"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.238"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.217"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.305"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.11"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.71%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +6.93%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.60%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.983"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.24%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05064"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7387"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.69%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.153"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02081"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.743"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.672"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.02%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.07"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8729"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.45%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.812"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4266"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9991"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "51.43"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +19.11%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.51"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.179"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.293"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.39%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.1210"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.84"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3941"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.19%  "
